$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 2938.103010863317
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 74457.80326820488

$ws.Range("B3").Value = 0.06328177979961902
$ws.Range("C3").Value = 1766.335244827366
$ws.Range("D3").Value = 157.8057217802531
$ws.Range("E3").Value = 5548678842208.939
$ws.Range("G3").Value = 5548678844133.144
